$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 45
# from serial date 45205 (2023-10-06) to 45206 (2023-10-07)
$ws.Range("C2:C45").Value = 45206
